# Adding notes for final deliverables.
# Inserts a new "Outline" slide as the first slide of the deck, and tidies
# up a run-split on the existing "A Digital Accelerator..." slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Insert the new Outline slide at position 1, using the "Title and
#    Content" layout (the same layout already used by the other slides).
# ---------------------------------------------------------------------
$outline = $p.Slides.Add(1, 2)

$title = $outline.Shapes.Item(1)
$titleTextRange = $title.TextFrame.TextRange
$titleTextRange.Text = "Outline"
$titleTextRange.LanguageID = "en-US"

$body = $outline.Shapes.Item(2)
$bodyTextRange = $body.TextFrame.TextRange
$bodyTextRange.Text = "an overview of the project and its objectives,`ra discussion of the architectural and timing design`rNew hardware architecture diagram`rexperiences and issues in implementation (the difficult parts)`rFixed point accuracy and precision`rFPGA resources`rDriver data movement`ra summary including lessons learned"
$bodyTextRange.LanguageID = "en-US"

# Indent the sub-bullets (IndentLevel is 1-based: level 2 == OOXML lvl="1")
$bodyTextRange.Paragraphs(3, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(5, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(6, 1).IndentLevel = 2
$bodyTextRange.Paragraphs(7, 1).IndentLevel = 2

# Give the body placeholder the same kind of name PowerPoint assigned in
# the authored deck.
$body.Name = "Content Placeholder 8"

# ---------------------------------------------------------------------
# 2. On the "A Digital Accelerator for Inverse Kinematics" slide, merge
#    the two runs that made up "...general algorithm on CPU" into one.
# ---------------------------------------------------------------------
$accelSlide = $p.Slides.Item(3)
$accelBody = $accelSlide.Shapes.Item(2).TextFrame.TextRange
$lastPara = $accelBody.Paragraphs($accelBody.Paragraphs().Count, 1)
$mergedTail = $lastPara.Characters(24, 52)
$mergedTail.Text = "s: compare against 10ms for general algorithm on CPU"
$mergedTail.LanguageID = "en-US"
